$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Year (A) and product_type (F) for the 10 new rows: these only ever reuse
# shared strings that already exist (or are plain numbers), so writing them
# up front does not disturb the new shared-string insertion order below.
$ws.Range("A12").Value = 2004
$ws.Range("F12").Value = "rulebook"
$ws.Range("A13").Value = 2004
$ws.Range("F13").Value = "supplement"
$ws.Range("A14").Value = 2005
$ws.Range("F14").Value = "replay"
$ws.Range("A15").Value = 2005
$ws.Range("F15").Value = "rulebook"
$ws.Range("A16").Value = 2006
$ws.Range("F16").Value = "supplement"
$ws.Range("A17").Value = 2006
$ws.Range("F17").Value = "supplement"
$ws.Range("A18").Value = 2007
$ws.Range("F18").Value = "supplement"
$ws.Range("A19").Value = 2007
$ws.Range("F19").Value = "supplement"
$ws.Range("A20").Value = 2008
$ws.Range("F20").Value = "supplement"
$ws.Range("A21").Value = 2010
$ws.Range("F21").Value = "supplement"

# japanese / english / publisher / image columns, written in the same order
# the original author entered them (matching the shared-string table order).
$ws.Range("B12").Value = "真・女神転生3‐NOCTURNE TRPG 東京受胎"
$ws.Range("C12").Value = "Shin Megami Tensei 3-NOCTURNE TRPG Tokyo conception"
$ws.Range("D12").Value = "Jive"
$ws.Range("E12").Value = "nocturne-rulebook.jpg"
$ws.Range("E13").Value = "nocturne-supplement.jpg"
$ws.Range("B13").Value = "アマラ深界 真・女神転生3‐NOCTURNE TRPGサプリメント"
$ws.Range("C13").Value = "Amara Deep World: Shin Megami Tensei 3-NOCTURNE TRPG Supplement"
$ws.Range("D13").Value = "Jive"
$ws.Range("B14").Value = "再会―See you again next world 真・女神転生3 NOCTURNE TRPGリプレイ"
$ws.Range("C14").Value = "Reunion―See you again next world: Shin Megami Tensei 3 NOCTURNE TRPG Replay"
$ws.Range("E14").Value = "nocturne-replay.jpg"
$ws.Range("D14").Value = "Jive"
$ws.Range("B15").Value = "真・女神転生TRPG　魔都東京200X"
$ws.Range("C15").Value = "Shin Megami Tensei TRPG Magic City Tokyo 200X"
$ws.Range("E15").Value = "200x-rulebook.jpg"
$ws.Range("D15").Value = "Jive"
$ws.Range("B16").Value = "金剛神界 真・女神転生TRPG魔都東京200Xサプリメント"
$ws.Range("E16").Value = "200x-magic-city-tokyo.jpg"
$ws.Range("C16").Value = "Kongo Shinkai: Shin Megami Tensei TRPG Magic City Tokyo 200X Supplement"
$ws.Range("D16").Value = "Jive"
$ws.Range("B17").Value = "TOKYOミレニアム 真・女神転生TRPG魔都東京200Xサプリメント"
$ws.Range("E17").Value = "200x-tokyo-millenium.jpg"
$ws.Range("C17").Value = "Tokyo Millennium: Shin Megami Tensei TRPG Magic City Tokyo 200X Supplement"
$ws.Range("D17").Value = "Jive"
$ws.Range("B18").Value = "セフィロトの魔界 真・女神転生TRPG魔都東京200Xサプリメント"
$ws.Range("C18").Value = "Sephiroto's Makai: Shin Megami Tensei TRPG Magic City Tokyo 200X Supplement"
$ws.Range("E18").Value = "200x-sephirots-hell.jpg"
$ws.Range("D18").Value = "Jive"
$ws.Range("B19").Value = "闇のプロファイル 真・女神転生TRPG魔都東京200X"
$ws.Range("E19").Value = "200x-dark-profile.jpg"
$ws.Range("C19").Value = "Profile of Darkness Shin Megami Tensei TRPG Magic City Tokyo 200X"
$ws.Range("D19").Value = "Jive"
$ws.Range("B21").Value = "異形科学 －ストレンジ・サイエンス 真・女神転生TRPG 魔都東京200X"
$ws.Range("C21").Value = "Variant Science-Strange Science Shin Megami Tensei TRPG Magic City Tokyo 200X"
$ws.Range("E20").Value = "200x-ragnarok.jpg"
$ws.Range("B20").Value = "ラグナロク 真・女神転生TRPG魔都東京200X"
$ws.Range("C20").Value = "Ragnarok: Shin Megami Tensei TRPG Magic City Tokyo 200X"
$ws.Range("D20").Value = "Jive"
$ws.Range("E21").Value = "200x-variant-science.jpg"
$ws.Range("D21").Value = "Jive"

# Column width adjustments for B (japanese) and C (english).
# The host's ColumnWidth setter quantizes to whole pixels (steps of 1/6
# character unit here), so the nearest representable widths to the
# authored 66.83203125 / 73.83203125 are reached via these inputs.
$ws.Columns.Item(2).ColumnWidth = 65.96
$ws.Columns.Item(3).ColumnWidth = 72.96

# Update selection to match the authored state
$ws.Range("E22").Select()
